# "Added cluster all feature" -- append a new "C4H2I2S" calibration block
# (8 rows: 21.0, 35.0, 50.0, 70.0, 100.0, 200.0, 300.0, 500.0) below the
# existing time-offset table, with a secondary offset column (D) and a
# black-font style applied to the label cells of the last 4 new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block: label, column B (offset us), column D (secondary offset us)
$rows = @(
    @{ Row = 94;  Label = "C4H2I2S_21.0";  B = -9; D = -9  },
    @{ Row = 95;  Label = "C4H2I2S_35.0";  B = 10; D = 10  },
    @{ Row = 96;  Label = "C4H2I2S_50.0";  B = 10; D = 10  },
    @{ Row = 97;  Label = "C4H2I2S_70.0";  B = 0;  D = 10  },
    @{ Row = 98;  Label = "C4H2I2S_100.0"; B = 0;  D = -8  },
    @{ Row = 99;  Label = "C4H2I2S_200.0"; B = 0;  D = -30 },
    @{ Row = 100; Label = "C4H2I2S_300.0"; B = 0;  D = -25 },
    @{ Row = 101; Label = "C4H2I2S_500.0"; B = 0;  D = -20 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Label
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}

# Rows 98-101 get a black-font style on the label cell (column A)
$ws.Range("A98:A101").Font.Color = 0

# Update the selection to reflect where the author ended up editing
$ws.Range("B98").Select()
